$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Metrics": update values and the selected/active cell
# ---------------------------------------------------------------------------
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 270148.56
$metrics.Range("B3").Value  = 232713.73
$metrics.Range("B4").Value  = 85064.42
$metrics.Range("B5").Value  = 10709
$metrics.Range("B6").Value  = 3665777.13
$metrics.Range("B7").Value  = 3111428.39
$metrics.Range("B8").Value  = 1048706.98
$metrics.Range("B9").Value  = 141397
$metrics.Range("B10").Value = 32131100.930999827
$metrics.Range("B11").Value = 19141298.460000001
$metrics.Range("B12").Value = 11330415.870000001
$metrics.Range("B13").Value = 1239024

# Update the sheet's selection to match the new active cell (C13)
$metrics.Activate() | Out-Null
$metrics.Range("C13").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "today": replace the formulas with their plain computed values
# ---------------------------------------------------------------------------
$today = $wb.Worksheets.Item("today")

$today.Range("B11").Value = 270148.56
$today.Range("B12").Value = 232713.73
$today.Range("B13").Value = 85064.42

$today.Range("B14").Value = 10709
$today.Range("B14").NumberFormat = "#,##0.00"

$today.Range("B15").Value = 3665777.13
$today.Range("B16").Value = 3111428.39
$today.Range("B17").Value = 1048706.98
$today.Range("B18").Value = 141397
$today.Range("B19").Value = 32131100.930999827
$today.Range("B20").Value = 19141298.460000001
$today.Range("B21").Value = 11330415.870000001
$today.Range("B22").Value = 1239024
